# Componenti da acquistare.xlsx update
# Adds a "(Boost?)" note to the sonar row, replaces the dual-option Stereo
# Camera block with a single ZED 2i entry, adds the StereoCamera cylinder /
# end-caps / GPS rows, and refreshes the totals at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- alignment / enum constants (standard Excel values) ----
$xlGeneral = 1
$xlLeft    = -4131
$xlCenter  = -4108
$xlBottom  = -4107

$eurFmt410 = "[`$€-410]\ #,##0.00;[RED]\-[`$€-410]\ #,##0.00"

# ---------------------------------------------------------------------
# Row 2 - ECHO Multibeam Imaging Sonar gains a "(Boost?)" second line
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "ECHO Multibeam Imaging Sonar`n(Boost?)"
$ws.Range("A2").WrapText = $true
$ws.Range("A2").HorizontalAlignment = $xlCenter
$ws.Range("A2").VerticalAlignment = $xlCenter
$ws.Rows.Item(2).RowHeight = 28.35

# ---------------------------------------------------------------------
# Row 5 - height shrinks now that the Stereo Camera block is single-option
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "Stereo Camera (da scegliere 2 o 4mm)"
$ws.Range("B5").Value = "ZED 2i (IP66)"
$ws.Range("C5").Value = "https://store.stereolabs.com/products/zed-2i?variant=41379929096348"
$ws.Range("D5").Value = "€499,00`n"
$ws.Range("E5").Value = "cella sottostante"
$ws.Range("F5").Value = "175.25 x 30.25 x 43.10 mm"
$ws.Range("G5").Value = "166 g"
$ws.Range("H5").Value = "1) 5 V (via USB)`n"
$ws.Range("I5").Value = "2 W"
$ws.Rows.Item(5).RowHeight = 28.35

# ---------------------------------------------------------------------
# Row 8 - Cilindro StereoCamera (new)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Cilindro StereoCamera"
$ws.Range("B8").Value = "BlueRobotics"
$ws.Hyperlinks.Add($ws.Range("C8"), "https://bluerobotics.com/store/watertight-enclosures/locking-series/wte-locking-tube-r1-vp/", "", "", "https://bluerobotics.com/store/watertight-enclosures/locking-series/wte-locking-tube-r1-vp/") | Out-Null
$ws.Range("C8").Formula = "https://bluerobotics.com/store/watertight-enclosures/locking-series/wte-locking-tube-r1-vp/`nAcrilico, 75 mm (diametro), 300mm (lunghezza)"
$ws.Range("D8").Value = "`$215.00=€197,71"
$ws.Range("G8").Value = "560 g"

$ws.Range("A8").WrapText = $true
$ws.Range("A8").HorizontalAlignment = $xlCenter
$ws.Range("A8").VerticalAlignment = $xlCenter
$ws.Range("C8").WrapText = $true
$ws.Range("C8").HorizontalAlignment = $xlCenter
$ws.Range("C8").VerticalAlignment = $xlCenter
$ws.Rows.Item(8).RowHeight = 28.35

# ---------------------------------------------------------------------
# Row 9 - Tappi Cilindro (new)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Tappi Cilindro"
$ws.Range("B9").Value = "BlueRobotics"
$ws.Hyperlinks.Add($ws.Range("C9"), "https://bluerobotics.com/store/watertight-enclosures/locking-series/wte-end-cap-vp/", "", "", "https://bluerobotics.com/store/watertight-enclosures/locking-series/wte-end-cap-vp/") | Out-Null
$ws.Range("C9").Formula = "https://bluerobotics.com/store/watertight-enclosures/locking-series/wte-end-cap-vp/ `nAlluminio, 1x 4 fori, 1x senza fori, 75mm"
$ws.Range("D9").Value = "`$28.00=€25,75"
$ws.Range("G9").Value = "97 g + 102 g"

$ws.Range("A9").WrapText = $true
$ws.Range("A9").HorizontalAlignment = $xlCenter
$ws.Range("A9").VerticalAlignment = $xlCenter
$ws.Range("C9").WrapText = $true
$ws.Range("C9").HorizontalAlignment = $xlCenter
$ws.Range("C9").VerticalAlignment = $xlCenter
$ws.Rows.Item(9).RowHeight = 28.35

# ---------------------------------------------------------------------
# Row 10 - GPS (da comprare?) (new)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "GPS (da comprare?)"
$ws.Range("B10").Value = "Ublox"
$ws.Hyperlinks.Add($ws.Range("C10"), "https://futuranet.it/prodotto/ublox-neo-6m-per-controller-di-volo-apm/", "", "", "https://futuranet.it/prodotto/ublox-neo-6m-per-controller-di-volo-apm/") | Out-Null
$ws.Range("C10").Formula = "https://futuranet.it/prodotto/ublox-neo-6m-per-controller-di-volo-apm/"
$ws.Range("D10").Value = 21
$ws.Range("D10").NumberFormat = $eurFmt410
$ws.Range("D10").HorizontalAlignment = $xlCenter
$ws.Range("D10").VerticalAlignment = $xlCenter
$ws.Range("E10").Value = "Cilindro superiore"
$ws.Range("F10").Value = "30x23x4 mm (modulo)`n25x25x8 mm (antenna)"
$ws.Range("F10").WrapText = $true
$ws.Range("F10").HorizontalAlignment = $xlCenter
$ws.Range("F10").VerticalAlignment = $xlCenter
$ws.Range("H10").Value = "3-5 V"
$ws.Rows.Item(10).RowHeight = 28.35

# ---------------------------------------------------------------------
# Row 12 - updated totals (row 11 intentionally left blank)
# ---------------------------------------------------------------------
$ws.Range("D12").Value = "Tot. € 9426,26"
$ws.Range("D12").HorizontalAlignment = $xlGeneral
$ws.Range("D12").VerticalAlignment = $xlCenter

$ws.Range("G12").Value = "Tot. 1868 g"
$ws.Range("G12").HorizontalAlignment = $xlLeft
$ws.Range("G12").VerticalAlignment = $xlBottom

$ws.Rows.Item(12).RowHeight = 13.8

# ---------------------------------------------------------------------
# Column C widens to fit the longer StereoCamera / end-cap descriptions
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 81

# ---------------------------------------------------------------------
# Restore the selection the author left the sheet on
# ---------------------------------------------------------------------
$ws.Range("G17:G18").Select() | Out-Null
